$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: insert the new "Workblocks" sheet right after "Settings" (so the
# tab order becomes Settings, Workblocks, Constants, Assets).
# ---------------------------------------------------------------------------
$settings   = $wb.Worksheets.Item(1)
$workblocks = $wb.Worksheets.Add($null, $settings)
$workblocks.Name = "Workblocks"

# Re-fetch sheet references AFTER the sheet collection changed.
$constants = $wb.Worksheets.Item("Constants")

# ---------------------------------------------------------------------------
# Step 2: capture the 14 "Workblock" rows that currently live in Constants
# (old rows 19-32) before they get removed from that sheet.
# ---------------------------------------------------------------------------
$taskRows = @()
for ($i = 19; $i -le 32; $i++) {
    $a = $constants.Cells.Item($i, 1).Value()
    $b = $constants.Cells.Item($i, 2).Value()
    $c = $constants.Cells.Item($i, 3).Value()
    $taskRows += , @($a, $b, $c)
}

# ---------------------------------------------------------------------------
# Step 3: populate the new Workblocks sheet.
# ---------------------------------------------------------------------------
$workblocks.Range("A1").Value = "Name"
$workblocks.Range("B1").Value = "Value"
$workblocks.Range("C1").Value = "Description"

# Re-use the existing bold-header style (same style already used by row 1 of
# the other sheets) by copying formats from an already-styled header cell,
# instead of re-deriving Bold/Size property-by-property (which leaves orphan
# cellXfs entries behind in this engine).
$settings.Range("A1").Copy() | Out-Null
$workblocks.Range("A1:C2").PasteSpecial(-4122) | Out-Null
$workblocks.Rows.Item(1).RowHeight = 18.75
$workblocks.Rows.Item(2).RowHeight = 18.75

for ($i = 0; $i -lt $taskRows.Count; $i++) {
    $r = 3 + $i
    $row = $taskRows[$i]
    $workblocks.Cells.Item($r, 1).Value = $row[0]
    $workblocks.Cells.Item($r, 2).Value = $row[1]
    $workblocks.Cells.Item($r, 3).Value = $row[2]
    if ($r -ne 13) {
        $workblocks.Cells.Item($r, 2).HorizontalAlignment = -4131
    }
}

$workblocks.Columns.Item(1).ColumnWidth = 41.140625
$workblocks.Columns.Item(2).ColumnWidth = 27.85546875
$workblocks.Columns.Item(3).ColumnWidth = 35.85546875

# ---------------------------------------------------------------------------
# Step 4: remove the Workblock rows from Constants and re-pad with blank
# rows (18-32), then shift everything from row 5 down by one (a blank
# separator row is inserted at row 5).
# ---------------------------------------------------------------------------
$constants.Range("A18:C32").Clear() | Out-Null
for ($i = 18; $i -le 32; $i++) {
    $constants.Rows.Item($i).RowHeight = 14.25
}

$constants.Rows.Item(5).Insert() | Out-Null
$constants.Rows.Item(5).RowHeight = 14.25

# ---------------------------------------------------------------------------
# Step 5: fix up tab-selection / active-cell state to match the target.
# ---------------------------------------------------------------------------
$constants.Range("A1:C1").Select() | Out-Null
$workblocks.Activate()
$workblocks.Range("B21").Select() | Out-Null

Write-Output "done"
